# The document is a diary: each entry is a "date" paragraph
# ("2022<ideographs>" split into 3 runs: "2" | "022" | "<rest>",
# the first/last runs carrying an eastAsia rFonts hint) followed by a
# plain single-run "description" paragraph (also eastAsia-hinted).
#
# The last paragraph in the document is an empty stub (just an empty
# <w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>) left
# over from a trailing Enter press. We need to turn it into a new
# "2022年6月7日星期二" date entry, plus a new description paragraph
# "晴，今天是高考第一天，上午考语文，下午考数学" right after it.
#
# Rather than editing that stub paragraph in place (which stubbornly
# keeps its empty <w:pPr> around, not matching the target), we insert
# two brand-new paragraphs after the last real content paragraph
# (which come out clean, with no <w:pPr>), fill those in, and then
# delete the original empty stub paragraph.

$d = $word.ActiveDocument

$lastContent = $d.Paragraphs($d.Paragraphs.Count - 1)

# Create two fresh paragraphs after the last real ("...赛龙舟。") entry.
$lastContent.Range.InsertParagraphAfter()
$datePara = $d.Paragraphs($lastContent.Index + 1)
$datePara.Range.InsertParagraphAfter()

$datePara = $d.Paragraphs($lastContent.Index + 1)
$descPara = $d.Paragraphs($lastContent.Index + 2)
$stubPara = $d.Paragraphs($lastContent.Index + 3)

# --- Build the date paragraph: "2" | "022" | "年6月7日星期二" ---
# Use the very first date paragraph in the doc as the formatting
# template: its 3 runs already show the hint / no-hint / hint pattern
# we need, and its 3rd run happens to be the same length (8 chars) as
# our target text, so a FormattedText copy carries the eastAsia hint
# over cleanly; we then patch the text back to the right digits.
$templateDate = $d.Paragraphs(1)
$tplStart = $templateDate.Range.Start

$dateStart = $datePara.Range.Start
$datePara.Range.Text = "2022年6月7日星期二"

$d.Range($dateStart, $dateStart + 1).FormattedText = `
    $d.Range($tplStart, $tplStart + 1).FormattedText
$d.Range($dateStart + 1, $dateStart + 4).FormattedText = `
    $d.Range($tplStart + 1, $tplStart + 4).FormattedText
$d.Range($dateStart + 4, $dateStart + 12).FormattedText = `
    $d.Range($tplStart + 4, $tplStart + 12).FormattedText
$d.Range($dateStart + 4, $dateStart + 12).Text = "年6月7日星期二"

# --- Build the description paragraph (single eastAsia-hinted run) ---
# Paragraph 4 ("中雨，今天是农历五月初四，明天就是端午节了。") is a
# single hinted run of exactly 22 characters, matching the length of
# our target description text, so it makes a perfect format donor.
$descPara = $d.Paragraphs($lastContent.Index + 2)
$templateDesc = $d.Paragraphs(4)
$tplDescStart = $templateDesc.Range.Start

$descText = "晴，今天是高考第一天，上午考语文，下午考数学"
$descStart = $descPara.Range.Start
$descPara.Range.Text = $descText

$d.Range($descStart, $descStart + 22).FormattedText = `
    $d.Range($tplDescStart, $tplDescStart + 22).FormattedText
$d.Range($descStart, $descStart + 22).Text = $descText

# --- Remove the original empty trailing stub paragraph ---
$descPara = $d.Paragraphs($lastContent.Index + 2)
$stubPara = $d.Paragraphs($lastContent.Index + 3)
$d.Range($descPara.Range.End - 1, $stubPara.Range.End).Delete()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
